# issue #5: stock data from json to db
# Add "category" column (after property_category) and trailing
# "source_file" / "index" columns to the 股票 (stock) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Insert new "category" column right after "property_category" ---
# (old column H keeps its data; new column I is inserted before the old I/J/K)
$ws.Columns("I").Insert()

# Header cell: copy the bold/bordered header style from H1, then set text
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "category"

# Data cells: copy the plain data style from H2:H6, then set text
$ws.Range("H2:H6").Copy() | Out-Null
$ws.Range("I2:I6").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"
$ws.Range("I4").Value = "normal"
$ws.Range("I5").Value = "normal"
$ws.Range("I6").Value = "normal"

# --- Append two new trailing columns: source_file (M) and index (N) ---
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:N1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$ws.Range("L2:L6").Copy() | Out-Null
$ws.Range("M2:M6").PasteSpecial(-4122) | Out-Null
$ws.Range("L2:L6").Copy() | Out-Null
$ws.Range("N2:N6").PasteSpecial(-4122) | Out-Null

$ws.Range("M2").Value = "tmpacad1"
$ws.Range("M3").Value = "tmpacad1"
$ws.Range("M4").Value = "tmpacad1"
$ws.Range("M5").Value = "tmpacad1"
$ws.Range("M6").Value = "tmpacad1"

$ws.Range("N2").Value = 72
$ws.Range("N3").Value = 73
$ws.Range("N4").Value = 74
$ws.Range("N5").Value = 75
$ws.Range("N6").Value = 76
